$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Cypher query text for the "CasesTab" row (B2): case query reordered so the
# WHERE clause comes after the MATCH clauses instead of right after the first MATCH.
$newB2 = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
 WHERE c.gender='MALE'
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# New Cypher query text for the stats query on the "CasesTab" row (C2): replaces the
# specimen-based Trials/Cases/Files count query with a file-based one.
$newC2 = @'

MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
    WHERE c.gender = "MALE"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$ws.Range("B2").Value = $newB2
$ws.Range("C2").Value = $newC2

# Row 2 shrank by one wrapped line once the query text changed; match the new
# auto-fitted row height.
$ws.Rows.Item(2).RowHeight = 195

# Move the active selection from B3 to B2.
$null = $ws.Range("B2").Select()
